$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pmclEmail = "sreekanth.pogula@senecaglobal.com"

$data = @(
    @("N1070", "Sreekanth Pogula", $pmclEmail, $pmclEmail, "sreekanth.pogula@senecaglobal.com"),
    @("N1071", "Bhargavi Gottumukkala", $pmclEmail, $pmclEmail, "bhargavi.gottumukkala@senecaglobal.com"),
    @("N1072", "Harnath Immani", $pmclEmail, $pmclEmail, "harnath.immani@senecaglobal.com"),
    @("N1073", "Jawahar Prudhivi", $pmclEmail, $pmclEmail, "jawahar.prudhivi@senecaglobal.com"),
    @("N1074", "Manisha Siram", $pmclEmail, $pmclEmail, "manisha.siram@senecaglobal.com"),
    @("N1079", "Nandini Yerrapothu", $pmclEmail, $pmclEmail, "nandini.yerrapothu@senecaglobal.com"),
    @("N1056", "Prakash Chandra", $pmclEmail, $pmclEmail, "prakash.chandra@senecaglobal.com"),
    @("N1000", "Ramesh Gidde", $pmclEmail, $pmclEmail, "ramesh.gidde@senecaglobal.com"),
    @("N1234", "Sabiha Sultana", $pmclEmail, $pmclEmail, "sabiha.sultana@senecaglobal.com"),
    @("N5363", "Shravani Deshpande", $pmclEmail, $pmclEmail, "shravani.deshpande@senecaglobal.com"),
    @("N3131", "Yagnabhargavi Penumacha", $pmclEmail, $pmclEmail, "yagnabhargavi.penumacha@senecaglobal.com")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
